$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The export window rolled forward by two days: the two oldest rows
# (2025-09-24 and 2025-09-25) drop off the front, every remaining row
# shifts up, and a new row for 2025-12-24 is appended at the end.
$ws.Rows("2:3").Delete()

# Append the new trailing row (now row 91) for 2025-12-24.
$lastRow = $ws.UsedRange.Rows.Count + 1

$dateCell = $ws.Cells.Item($lastRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-12-24"
$dateCell.ClearFormats()

$ws.Cells.Item($lastRow, 2).Value = 0
$ws.Cells.Item($lastRow, 3).Value = 0
